# Weekly data refresh: a new price observation is inserted as row 52
# (the sheet's data rows are sorted newest-first after the header/blank
# rows), pushing all the former rows 52-75 down by one to become 53-76.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 52; everything below (old 52..75)
# shifts down to 53..76, carrying its formatting and values with it.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44627
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = "Poroto verde"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 60
$ws.Range("K52").Value = 29000
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = 29500
$ws.Range("N52").Value = "$/saco 25 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1180
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"

# Match the date column's number format style used by the rest of the
# column (style index 2 / numFmtId 165 "YYYY-MM-DD HH:MM:SS").
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat()
